$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# D17 and D18 were stored as text (inline strings); convert them to real numbers.
$ws.Range("D17").Value = 500480
$ws.Range("D18").Value = 532478

# Append the new row 19 (COLPAL) at the bottom of the table.
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "COLPAL"
$ws.Range("C19").Value = "Colgate Palmolive (india) Limited"
# bsecode kept as text for this row (matches source data), so force text format
# before writing the numeric-looking string.
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "500830"
$ws.Range("E19").Value = 0.5600000000000001
$ws.Range("F19").Value = 2952.6
$ws.Range("G19").Value = 251795
$ws.Range("H19").Value = "day"
$ws.Range("I19").Value = "14/06/2024 10:32:27"
